$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New experiment name used for the Rambam 12 extract rows
$expName = "Rambam 12"
$expDate = 43503   # 2019-02-07, serial date number
$expNum  = 12

$basePath = "W:\phkinnerets\storage\analysis\Niv\rambam12 extract\test extract\2019_02_07\Capture "

for ($i = 1; $i -le 8; $i++) {
    $row = 7 + $i

    # Copy formatting from the analogous cells in row 2 so the new
    # rows pick up the same styles (date format, centered alignment)
    # that the rest of the table uses.
    $ws.Range("A2").Copy()
    $ws.Range("A$row").PasteSpecial(-4122)
    $ws.Range("A$row").Value = $expDate

    $ws.Range("C2").Copy()
    $ws.Range("B$row").PasteSpecial(-4122)
    $ws.Range("C$row").PasteSpecial(-4122)

    $ws.Range("B$row").Value = $expName
    $ws.Range("C$row").Value = $expNum
    $ws.Range("D$row").Value = "$basePath$i\"
}

[void]$ws.Range("B19").Select()
